# Update the "Förändrad" (Changed) date column C for rows 2-6
# from serial date 45175 (2023-09-06) to serial date 45183 (2023-09-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45183
}
